# Auto-generated-assisted edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H48").Value = 4533.3335
$ws.Range("J48").Value = 8333.333000000001
$ws.Range("L48").Value = 24999.999
$ws.Range("N48").Value = -25583.999

$ws.Range("H56").Value = 4533.3335
$ws.Range("J56").Value = 8333.333000000001
$ws.Range("L56").Value = 24999.999
$ws.Range("N56").Value = -26067.999

$ws.Range("H70").Value = 3600451.5
$ws.Range("I70").Value = 5994419.5
$ws.Range("J70").Value = 9500
$ws.Range("K70").Value = 17983258.5
$ws.Range("L70").Value = 28500
$ws.Range("M70").Value = -17982988.5
$ws.Range("N70").Value = -29040

$ws.Range("H73").Value = 3600451.5
$ws.Range("I73").Value = 5994419.5
$ws.Range("J73").Value = 9500
$ws.Range("K73").Value = 17983258.5
$ws.Range("L73").Value = 28500
$ws.Range("M73").Value = -17982322.5
$ws.Range("N73").Value = -30372

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 1323.7693
$ws.Range("I107").Value = 1340.4546
$ws.Range("J107").Value = 1232
$ws.Range("K107").Value = 1340.4546
$ws.Range("L107").Value = 1232
$ws.Range("M107").Value = 579.5454
$ws.Range("N107").Value = -5072

$ws.Range("H137").Value = 15627361
$ws.Range("J137").Value = 2575.558
$ws.Range("L137").Value = 7726.674
$ws.Range("N137").Value = -12826.674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2124.5
$ws.Range("I61").Value = 1015.1579
$ws.Range("J61").Value = 4466.4443
$ws.Range("K61").Value = 1015.1579
$ws.Range("L61").Value = 4466.4443
$ws.Range("M61").Value = -803.1579
$ws.Range("N61").Value = -4890.4443

$ws.Range("H97").Value = 265.25806
$ws.Range("I97").Value = 118.46154
$ws.Range("K97").Value = 118.46154
$ws.Range("M97").Value = 377.53846

$ws.Range("H122").Value = 4592.684
$ws.Range("I122").Value = 4051.75
$ws.Range("K122").Value = 12155.25
$ws.Range("M122").Value = -9705.25

$ws.Range("H132").Value = 364679.06
$ws.Range("I132").Value = 197106.84
$ws.Range("K132").Value = 591320.52
$ws.Range("M132").Value = -588790.52

$ws.Range("H136").Value = 2124.5
$ws.Range("I136").Value = 1015.1579
$ws.Range("J136").Value = 4466.4443
$ws.Range("K136").Value = 3045.4737
$ws.Range("L136").Value = 13399.3329
$ws.Range("M136").Value = -495.4737
$ws.Range("N136").Value = -18499.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 41670772
$ws.Range("I86").Value = 55559788
$ws.Range("K86").Value = 55559788
$ws.Range("M86").Value = -55558665

$ws.Range("H89").Value = 41670772
$ws.Range("I89").Value = 55559788
$ws.Range("K89").Value = 277798940
$ws.Range("M89").Value = -277793324

$ws.Range("H97").Value = 7557.25
$ws.Range("I97").Value = 7557.25
$ws.Range("K97").Value = 7557.25
$ws.Range("M97").Value = -6566.25

$ws.Range("H107").Value = 13525641
$ws.Range("I107").Value = 8165
$ws.Range("J107").Value = 55580010
$ws.Range("K107").Value = 8165
$ws.Range("L107").Value = 55580010
$ws.Range("M107").Value = -6245
$ws.Range("N107").Value = -55583850

$ws.Range("H134").Value = 2130.7917
$ws.Range("I134").Value = 1642.7273
$ws.Range("J134").Value = 7499.5
$ws.Range("K134").Value = 4928.1819
$ws.Range("L134").Value = 22498.5
$ws.Range("M134").Value = -2393.1819
$ws.Range("N134").Value = -27568.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3202.1738
$ws.Range("I31").Value = 2920.0715
$ws.Range("K31").Value = 2920.0715
$ws.Range("M31").Value = -2625.0715

$ws.Range("H34").Value = 3202.1738
$ws.Range("I34").Value = 2920.0715
$ws.Range("K34").Value = 2920.0715
$ws.Range("M34").Value = -2718.0715

$ws.Range("H58").Value = 1980.8422
$ws.Range("I58").Value = 1563
$ws.Range("J58").Value = 2130.0715
$ws.Range("K58").Value = 1563
$ws.Range("L58").Value = 2130.0715
$ws.Range("M58").Value = -1360
$ws.Range("N58").Value = -2536.0715

$ws.Range("H107").Value = 1648.3334
$ws.Range("I107").Value = 519.6667
$ws.Range("J107").Value = 2212.6667
$ws.Range("K107").Value = 519.6667
$ws.Range("L107").Value = 2212.6667
$ws.Range("M107").Value = 1400.3333
$ws.Range("N107").Value = -6052.6667

$ws.Range("H134").Value = 2670.8276
$ws.Range("J134").Value = 3507.1538
$ws.Range("L134").Value = 10521.4614
$ws.Range("N134").Value = -15591.4614

$ws.Range("H136").Value = 1980.8422
$ws.Range("I136").Value = 1563
$ws.Range("J136").Value = 2130.0715
$ws.Range("K136").Value = 4689
$ws.Range("L136").Value = 6390.2145
$ws.Range("M136").Value = -2139
$ws.Range("N136").Value = -11490.2145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 150.66667
$ws.Range("I2").Value = 125
$ws.Range("J2").Value = 158
$ws.Range("K2").Value = 750
$ws.Range("L2").Value = 948
$ws.Range("M2").Value = -637
$ws.Range("N2").Value = -1174

$ws.Range("H38").Value = 124.72727
$ws.Range("I38").Value = 146.88889
$ws.Range("J38").Value = 25
$ws.Range("K38").Value = 440.66667
$ws.Range("L38").Value = 75
$ws.Range("M38").Value = -93.66667000000001
$ws.Range("N38").Value = -769

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10309.5
$ws.Range("I70").Value = 10312.846
$ws.Range("J70").Value = 10300.8
$ws.Range("K70").Value = 10312.846
$ws.Range("L70").Value = 10300.8
$ws.Range("M70").Value = -10042.846
$ws.Range("N70").Value = -10840.8

$ws.Range("H73").Value = 10309.5
$ws.Range("I73").Value = 10312.846
$ws.Range("J73").Value = 10300.8
$ws.Range("K73").Value = 10312.846
$ws.Range("L73").Value = 10300.8
$ws.Range("M73").Value = -9376.846
$ws.Range("N73").Value = -12172.8

$ws.Range("H80").Value = 2169.75
$ws.Range("J80").Value = 3499.6667
$ws.Range("L80").Value = 3499.6667
$ws.Range("N80").Value = -5495.6667

$ws.Range("H83").Value = 2169.75
$ws.Range("J83").Value = 3499.6667
$ws.Range("L83").Value = 17498.3335
$ws.Range("N83").Value = -27482.3335

$ws.Range("H122").Value = 3346.0908
$ws.Range("I122").Value = 2850
$ws.Range("J122").Value = 3941.4
$ws.Range("K122").Value = 8550
$ws.Range("L122").Value = 11824.2
$ws.Range("M122").Value = -6100
$ws.Range("N122").Value = -16724.2

$ws.Range("H132").Value = 336719.6
$ws.Range("I132").Value = 419649.6
$ws.Range("K132").Value = 1258948.8
$ws.Range("M132").Value = -1256418.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2601.9614
$ws.Range("I16").Value = 2147.1428
$ws.Range("K16").Value = 2147.1428
$ws.Range("M16").Value = -1977.1428

$ws.Range("H55").Value = 364.31818
$ws.Range("I55").Value = 399.75
$ws.Range("K55").Value = 399.75
$ws.Range("M55").Value = -226.75

$ws.Range("H61").Value = 6925.8276
$ws.Range("I61").Value = 8081.364
$ws.Range("K61").Value = 8081.364
$ws.Range("M61").Value = -7879.364

$ws.Range("H113").Value = 6925.8276
$ws.Range("I113").Value = 8081.364
$ws.Range("K113").Value = 8081.364
$ws.Range("M113").Value = -5911.364

$ws.Range("H136").Value = 5185.1904
$ws.Range("I136").Value = 1469.2307
$ws.Range("K136").Value = 4407.6921
$ws.Range("M136").Value = -1857.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9000
$ws.Range("I54").Value = 9000
$ws.Range("K54").Value = 9000
$ws.Range("M54").Value = -8480

$ws.Range("H107").Value = 456.07144
$ws.Range("I107").Value = 443.81818
$ws.Range("K107").Value = 1331.45454
$ws.Range("M107").Value = 588.54546

$ws.Range("H113").Value = 2961.9546
$ws.Range("I113").Value = 1396.8235
$ws.Range("J113").Value = 8283.4
$ws.Range("K113").Value = 4190.470499999999
$ws.Range("L113").Value = 24850.2
$ws.Range("M113").Value = -2020.470499999999
$ws.Range("N113").Value = -29190.2

$ws.Range("H122").Value = 3074.3125
$ws.Range("I122").Value = 1753
$ws.Range("K122").Value = 5259
$ws.Range("M122").Value = -2809

$ws.Range("H136").Value = 4804.4614
$ws.Range("I136").Value = 3641.5789
$ws.Range("J136").Value = 7960.857
$ws.Range("K136").Value = 10924.7367
$ws.Range("L136").Value = 23882.571
$ws.Range("M136").Value = -8374.736699999999
$ws.Range("N136").Value = -28982.571
